$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "D" = 0.5579999999999999
    "E" = 0.762
    "K" = 38.7
    "L" = 0.1454887218045113
    "M" = 103
    "N" = 0.1815937940761636
    "O" = 2.661498708010336
    "P" = 35.4
    "Q" = 0.06241184767277855
    "R" = 0.9147286821705425
    "S" = 67.59999999999999
    "T" = 0.6563106796116505
    "U" = 435.7
    "V" = 0.7681593794076162
    "W" = 0.03608391608391609
    "X" = 0.2885869086583041
    "Y" = -0.252502992574388
    "Z" = 0.08027765202957597
    "AB" = 0.1430344091953701
    "AC" = -0.1430344091953701
    "AD" = 3212.1
    "AF" = 3212.1
    "AG" = 2776.4
    "AH" = 0.8499192972243537
    "AI" = 0.7265222111643898
    "AJ" = 0.8303624835506639
    "AK" = 0.6966252665913938
}

foreach ($row in @(2, 3)) {
    foreach ($col in $values.Keys) {
        $ws.Range("$col$row").Value = $values[$col]
    }
}
